# Generate Report for Handoff
# A new source file (96781129-9ea5-4218-8114-ec6b2e9ac410.md) was handed off.
# Insert a new row 2 in every sheet (Overview, zh-cn, de-de) for the new file,
# pushing the existing e5462a4d... row down to row 3, and refresh the tables.

$wb = $excel.ActiveWorkbook

$ghBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1c0e6ce0b157006395c17eb086d609488af6a033/e2e/"
$newFile = "96781129-9ea5-4218-8114-ec6b2e9ac410.md"
$oldFile = "e5462a4d-dc0a-4c4b-844b-8d2791006370.md"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")

# Push the existing data row down from row 2 to row 3 (keeps cell types/styles).
$wsOv.Rows(2).Insert()

# Re-point the hyperlink that the insert left dangling on B2 so it sits on B3.
foreach ($h in @($wsOv.Hyperlinks)) { $h.Delete() }

# New row (row 2): the newly handed-off file.
$wsOv.Range("A2").Value = $newFile
$wsOv.Range("C2").Value = ".md"
$wsOv.Range("D2").Value = ""
$wsOv.Range("E2").Value = "Ready for handoff"
$wsOv.Range("F2").Value = "Ready for handoff"
$wsOv.Range("G2").Value = "2016-08-23 12:39:50"
$wsOv.Hyperlinks.Add($wsOv.Range("B2"), ($ghBase + $newFile), "", "", ("e2e\" + $newFile)) | Out-Null

# Row 3: restore the hyperlink for the file that was shifted down.
$wsOv.Hyperlinks.Add($wsOv.Range("B3"), ($ghBase + $oldFile), "", "", ("e2e\" + $oldFile)) | Out-Null

# Grow the Overview table / autofilter to include the new row.
$loOv = $wsOv.ListObjects.Item(1)
$loOv.Resize($wsOv.Range("A1:G3"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Rows(2).Insert()
foreach ($h in @($wsZh.Hyperlinks)) { $h.Delete() }

$wsZh.Range("A2").Value = $newFile
$wsZh.Range("B2").Value = ".md"
$wsZh.Range("C2").Value = "Ready for handoff"
$wsZh.Range("D2").Value = "e2e"
$wsZh.Range("E2").Value = "ht"
$wsZh.Range("F2").Value = "False"
$wsZh.Range("G2").Value = "96781129-9ea5-4218-8114-ec6b2e9ac410.1ae7c19e9c37af65c85f62d46bc551ea10f48df6.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-23 12:39:45"
$wsZh.Range("I2").Value = ""
$wsZh.Range("J2").Value = ""
$wsZh.Range("K2").Value = "0001-01-01 00:00:00"
$wsZh.Range("L2").Value = ""
$wsZh.Range("M2").Value = "True"
$wsZh.Range("N2").Value = ""
$wsZh.Range("O2").Value = "False"
$wsZh.Range("P2").Value = ""
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), ($ghBase + $newFile), "", "", $newFile) | Out-Null

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), ($ghBase + $oldFile), "", "", $oldFile) | Out-Null

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P3"))

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Rows(2).Insert()
foreach ($h in @($wsDe.Hyperlinks)) { $h.Delete() }

$wsDe.Range("A2").Value = $newFile
$wsDe.Range("B2").Value = ".md"
$wsDe.Range("C2").Value = "Ready for handoff"
$wsDe.Range("D2").Value = "e2e"
$wsDe.Range("E2").Value = "ht"
$wsDe.Range("F2").Value = "False"
$wsDe.Range("G2").Value = "96781129-9ea5-4218-8114-ec6b2e9ac410.1ae7c19e9c37af65c85f62d46bc551ea10f48df6.de-de.xlf"
$wsDe.Range("H2").Value = "2016-08-23 12:39:50"
$wsDe.Range("I2").Value = ""
$wsDe.Range("J2").Value = ""
$wsDe.Range("K2").Value = "0001-01-01 00:00:00"
$wsDe.Range("L2").Value = ""
$wsDe.Range("M2").Value = "True"
$wsDe.Range("N2").Value = ""
$wsDe.Range("O2").Value = "False"
$wsDe.Range("P2").Value = ""
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), ($ghBase + $newFile), "", "", $newFile) | Out-Null

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), ($ghBase + $oldFile), "", "", $oldFile) | Out-Null

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P3"))

Write-Host "Report regenerated for handoff of $newFile"
